# V2.1_Beta released: --update stopwords corpus
# Rebuild the DM hotwords table (column B = word, column C = frequency share)
# against the refreshed stopwords corpus. Some previously-listed words are
# now filtered as stopwords (熟悉, 相关, 用户, 工作, 负责) and drop out of the
# table, while newly-surfaced words (数学, 要求, 平台, Spark, 任职) take their
# place at the tail of the ranking. All frequency shares are recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "数据挖掘"
$ws.Range("C2").Value  = 0.261756970541316
$ws.Range("B3").Value  = "算法"
$ws.Range("C3").Value  = 0.174339580426122
$ws.Range("B4").Value  = "机器学习"
$ws.Range("C4").Value  = 0.152839588654459
$ws.Range("B5").Value  = "数据"
$ws.Range("C5").Value  = 0.129323562769115
$ws.Range("B6").Value  = "经验"
$ws.Range("C6").Value  = 0.122199616806708
$ws.Range("B7").Value  = "优先"
$ws.Range("C7").Value  = 0.105768334142884
$ws.Range("B8").Value  = "分析"
$ws.Range("C8").Value  = 0.0744094388538612
$ws.Range("B9").Value  = "能力"
$ws.Range("C9").Value  = 0.0725960352426972
$ws.Range("B10").Value = "数据分析"
$ws.Range("C10").Value = 0.070042580339844
$ws.Range("B11").Value = "模型"
$ws.Range("C11").Value = 0.0694139430530066
$ws.Range("B12").Value = "建模"
$ws.Range("C12").Value = 0.0661276099031103
$ws.Range("B13").Value = "挖掘"
$ws.Range("C13").Value = 0.0652240897031835
$ws.Range("B14").Value = "Python"
$ws.Range("C14").Value = 0.0599854299557824
$ws.Range("B15").Value = "以上学历"
$ws.Range("C15").Value = 0.0573462361267485
$ws.Range("B16").Value = "岗位职责"
$ws.Range("C16").Value = 0.0540339004938619
$ws.Range("B17").Value = "开发"
$ws.Range("C17").Value = 0.0531218941676708
$ws.Range("B18").Value = "海量"
$ws.Range("C18").Value = 0.0515685664314385
$ws.Range("B19").Value = "业务"
$ws.Range("C19").Value = 0.0487920824017798
$ws.Range("B20").Value = "技术"
$ws.Range("C20").Value = 0.0440095558918557
$ws.Range("B21").Value = "Hadoop"
$ws.Range("C21").Value = 0.0425923943293341
$ws.Range("B22").Value = "团队"
$ws.Range("C22").Value = 0.0408008927865718
$ws.Range("B23").Value = "Java"
$ws.Range("C23").Value = 0.040127239673617
$ws.Range("B24").Value = "优化"
$ws.Range("C24").Value = 0.0393481806825429
$ws.Range("B25").Value = "应用"
$ws.Range("C25").Value = 0.0382079466474276
$ws.Range("B26").Value = "熟练掌握"
$ws.Range("C26").Value = 0.0376833239386981
$ws.Range("B27").Value = "数学"
$ws.Range("C27").Value = 0.036262559575537
$ws.Range("B28").Value = "要求"
$ws.Range("C28").Value = 0.0361490326800581
$ws.Range("B29").Value = "平台"
$ws.Range("C29").Value = 0.0361250523559445
$ws.Range("B30").Value = "Spark"
$ws.Range("C30").Value = 0.0350599773257541
$ws.Range("B31").Value = "任职"
$ws.Range("C31").Value = 0.0345023314386803
